$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Vstart/Vstop measured values
$ws.Range("B2").Value = 12.6
$ws.Range("B3").Value = 12.2

# Update Ren1/Ren2 resistor values
$ws.Range("D4").Value = 150000
$ws.Range("D5").Value = 15000

# New helper formulas in column I
$ws.Range("I3").Formula = "=12.2*(D5/(D5+D4))"
$ws.Range("I4").Formula = "=12.6*D5/(D5+D4)"

# Update selection to match the saved view
$ws.Range("D6").Select()
